# Applies the cryptos.xlsx price/volume refresh described by the commit
# "Updated cryptos list on Fri May 19 22:59:05 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.916.68"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "1.816.06"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.13%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "309.54"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("E6").Value = "  +0.13%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4686"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.33%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3697"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.47%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07374"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.44%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.8711"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D12").Value = "1.847.20"
$ws.Range("E12").Value = "  +4.14%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.383"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.48%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.526"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.34%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.07075"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.28%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "91.75"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("E17").Value = "  +0.12%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008720"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("E19").Value = "  +0.04%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.75"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "26.950.46"
$ws.Range("E21").Value = "  +0.14%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.321"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("E23").Value = "  -2.65%  "
$ws.Range("D24").Value = "2.053.63"
$ws.Range("E24").Value = "  +2.60%  "
$ws.Range("E25").Value = "  -1.59%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "150.75"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.53%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.180"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.34%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "18.35"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("E29").Value = "  +0.79%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "116.28"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.79%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.08963"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.79%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.7695"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.30%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.166"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("E34").Value = "  +0.68%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.912"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("E36").Value = "  +0.15%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.087"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.90%  "
$ws.Range("E38").Value = "  +1.14%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01964"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.39%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.941"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("E41").Value = "  +0.39%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.5339"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.36%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.351"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -3.52%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.1658"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.24%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "8.454"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("E46").Value = "  -2.86%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "10.50"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.674"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.19%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "103.54"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.62%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.06304"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.37%  "
